$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Master")

# Data for the four new rows to append (rows 217-220)
$rows = @(
    @{
        B = "Physics-based full-band GaN high-electron-mobility transistor simulation suggests upper bound of LO phonon lifetime"
        C = 2026
        D = "AIP Publishing"
        E = "Journal of Applied Physics"
        F = "Dastider, Ankan Ghosh; Grupen, Matt; Tunga, Ashwin; Rakheja, Shaloo"
        H = "10.1063/5.0315424"
        I = "https://doi.org/10.1063/5.0315424"
        J = "Journal"
        K = "n-FET"
        L = "Experiment"
        M = "Transport"
        Q = "Physics-based full-band GaN high-electron-mobility transistor simulation suggests upper bound of LO phonon lifetime"
        R = "High"
        S = "2026-02-18"
    },
    @{
        B = "Novel FOM-enhanced sided-shield gate trench MOSFET with super-junction structure"
        C = 2026
        D = "IOP Publishing"
        E = "Semiconductor Science and Technology"
        F = "Yu, Hincheung; Sun, Yabin; Li, Xiaojin; Shi, Yanling; Shen, Yang; Ye, Bingyi; Zhang, Yuhang; Liu, Ziyu"
        H = "10.1088/1361-6641/ae4777"
        I = "https://doi.org/10.1088/1361-6641/ae4777"
        J = "Journal"
        K = "n-FET"
        L = "TCAD"
        M = "Gate Stack"
        Q = "Novel FOM-enhanced sided-shield gate trench MOSFET with super-junction structure"
        R = "High"
        S = "2026-02-18"
    },
    @{
        B = "Novel FOM-enhanced sided-shield gate trench MOSFET with super-junction structure"
        C = 2026
        D = "IOP Publishing"
        E = "Semiconductor Science and Technology"
        F = "Yu, Hincheung; Sun, Yabin; Li, Xiaojin; Shi, Yanling; Shen, Yang; Ye, Bingyi; Zhang, Yuhang; Liu, Ziyu"
        H = "10.1088/1361-6641/ae4777"
        I = "https://doi.org/10.1088/1361-6641/ae4777"
        J = "Journal"
        K = "n-FET"
        L = "TCAD"
        M = "Gate Stack"
        Q = "Novel FOM-enhanced sided-shield gate trench MOSFET with super-junction structure"
        R = "High"
        S = "2026-02-18"
    },
    @{
        B = "Physics-based full-band GaN high-electron-mobility transistor simulation suggests upper bound of LO phonon lifetime"
        C = 2026
        D = "AIP Publishing"
        E = "Journal of Applied Physics"
        F = "Dastider, Ankan Ghosh; Grupen, Matt; Tunga, Ashwin; Rakheja, Shaloo"
        H = "10.1063/5.0315424"
        I = "https://doi.org/10.1063/5.0315424"
        J = "Journal"
        K = "n-FET"
        L = "Experiment"
        M = "Transport"
        Q = "Physics-based full-band GaN high-electron-mobility transistor simulation suggests upper bound of LO phonon lifetime"
        R = "High"
        S = "2026-02-18"
    }
)

$startRow = 217
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = ""
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = ""
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 9).Value = $data.I
    $ws.Cells.Item($r, 10).Value = $data.J
    $ws.Cells.Item($r, 11).Value = $data.K
    $ws.Cells.Item($r, 12).Value = $data.L
    $ws.Cells.Item($r, 13).Value = $data.M
    $ws.Cells.Item($r, 14).Value = ""
    $ws.Cells.Item($r, 15).Value = ""
    $ws.Cells.Item($r, 16).Value = ""
    $ws.Cells.Item($r, 17).Value = $data.Q
    $ws.Cells.Item($r, 18).Value = $data.R
    $ws.Cells.Item($r, 19).Value = $data.S
    $ws.Cells.Item($r, 20).Value = ""
}
